$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("M6").Value = "Thomas G. Donlon"
$ws.Range("A8").Value = "Volume 31   Number  39"
$ws.Range("C9").Value = "Report Covering the Week  9/23/2024  Through  9/29/2024"

# --- Column E width (matches bestFit recompute after data change) ---
$ws.Columns("E").ColumnWidth = 6.168446

# Row 15
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50

# Row 16
$ws.Range("C16").Value = "'0"
$ws.Range("D16").Value = "'0"
$ws.Range("E16").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 16.666666666666
$ws.Range("M16").Value = -37.313432835820

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("I17").Value = 162
$ws.Range("J17").Value = 158
$ws.Range("K17").Value = 2.531645569620
$ws.Range("L17").Value = -4.142011834319
$ws.Range("M17").Value = 37.288135593220

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 71.428571428571
$ws.Range("I18").Value = 91
$ws.Range("J18").Value = 87
$ws.Range("K18").Value = 4.597701149425
$ws.Range("L18").Value = -35.460992907801
$ws.Range("M18").Value = -67.957746478873

# Row 19
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 20
$ws.Range("F19").Value = 51
$ws.Range("H19").Value = 6.25
$ws.Range("I19").Value = 407
$ws.Range("J19").Value = 428
$ws.Range("K19").Value = -4.906542056074
$ws.Range("L19").Value = -12.473118279569
$ws.Range("M19").Value = 40.344827586206

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 144
$ws.Range("J20").Value = 118
$ws.Range("K20").Value = 22.033898305084
$ws.Range("L20").Value = 67.441860465116
$ws.Range("M20").Value = 34.579439252336

# Row 21
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 14.285714285714
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = 21.739130434782
$ws.Range("I21").Value = 905
$ws.Range("J21").Value = 871
$ws.Range("K21").Value = 3.903559127439
$ws.Range("L21").Value = -4.029692470837
$ws.Range("M21").Value = -3.723404255319
$ws.Range("N21").Value = -76.830517153097

# Row 22
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Value = "'0"
$ws.Range("E22").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("L22").Value = 300

# Row 24
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 76
$ws.Range("G24").Value = 69
$ws.Range("H24").Value = 10.144927536231
$ws.Range("I24").Value = 734
$ws.Range("J24").Value = 815
$ws.Range("K24").Value = -9.938650306748
$ws.Range("L24").Value = -13.647058823529
$ws.Range("M24").Value = 5.763688760806

# Row 25
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 6.666666666666
$ws.Range("I25").Value = 180
$ws.Range("J25").Value = 239
$ws.Range("K25").Value = -24.686192468619
$ws.Range("L25").Value = -29.411764705882

# Row 26
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 62.5
$ws.Range("F26").Value = 39
$ws.Range("H26").Value = 34.482758620689
$ws.Range("I26").Value = 318
$ws.Range("J26").Value = 261
$ws.Range("K26").Value = 21.839080459770
$ws.Range("L26").Value = 23.735408560311
$ws.Range("M26").Value = -1.547987616099

# Row 27
$ws.Range("C27").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 500
$ws.Range("I28").Value = 52
$ws.Range("K28").Value = 6.122448979591
$ws.Range("L28").Value = -14.754098360655

# Row 33
$ws.Range("D33").Value = "'0"
$ws.Range("E33").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E33").PasteSpecial(-4122)
